$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the source inlineStr cells).
$ws.Range("D2").Value = "25.932.68"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "1.637.45"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.865.19"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.670.80"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "25.938.93"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "1.138.02"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "1.775.06"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -0.89%  "
